$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 322356.9755811149
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.2204408817635271
$ws.Range("F2").Value = 0.361247947454844
$ws.Range("G2").Value = 0.000246542046427866

# Row 3
$ws.Range("C3").Value = 51617607.23576894
$ws.Range("D3").Value = 0.394096529716793
$ws.Range("E3").Value = 0.3950946414289523
$ws.Range("F3").Value = 0.3945949544032483
$ws.Range("G3").Value = 0.0000574948900384189

# Row 4
$ws.Range("C4").Value = 25826015.21049996
$ws.Range("E4").Value = 0.3737669954678752
$ws.Range("F4").Value = 0.5441490393945274
$ws.Range("G4").Value = 0.00005907972616084886

# Row 5
$ws.Range("C5").Value = 51596401.59112387
$ws.Range("D5").Value = 0.7902796271637816
$ws.Range("E5").Value = 0.7911223673687017
$ws.Range("F5").Value = 0.7907007727151613
$ws.Range("G5").Value = 0.0001150887003738726

# Row 6
$ws.Range("C6").Value = 40849222.90134574
$ws.Range("D6").Value = 0.394096529716793
$ws.Range("E6").Value = 0.3950946414289523
$ws.Range("F6").Value = 0.3945949544032483
$ws.Range("G6").Value = 0.00007265128786499048

# Row 7
$ws.Range("C7").Value = 25846466.80855933
$ws.Range("E7").Value = 0.7709944014929352
$ws.Range("F7").Value = 0.8706909528827337
$ws.Range("G7").Value = 0.0001948458374900968

# Row 8
$ws.Range("C8").Value = 40832445.7539193
$ws.Range("D8").Value = 0.7902796271637816
$ws.Range("E8").Value = 0.7911223673687017
$ws.Range("F8").Value = 0.7907007727151613
$ws.Range("G8").Value = 0.000145427556283985

# Row 9
$ws.Range("C9").Value = 20443849.61398738
$ws.Range("E9").Value = 0.3737669954678752
$ws.Range("F9").Value = 0.5441490393945274
$ws.Range("G9").Value = 0.0000746333951418978

# Row 10
$ws.Range("C10").Value = 20460030.28549328
$ws.Range("E10").Value = 0.7709944014929352
$ws.Range("F10").Value = 0.8706909528827337
$ws.Range("G10").Value = 0.0002461421806909273

